# add results from latest run
#
# Refreshes the nowcast figures for the existing dates (rows 2-11) and
# appends a new row for the latest run date (2025-08-30), mirroring how
# this table is regenerated each time new survey/production data comes in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Refresh columns B:K for the already-present rows (2-11) ---------------

$rowData = @{
    2 = @(0.3083194854526926, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    3 = @(0.31486425940396073, 0, 0.008124571327022257, -0.00013448109970777696, -0.00015932350588940234, -0.0003105033614556734, -0.0000019455427035020168, -0.0010004638109840493, 0, 0.00002691994498627004)
    4 = @(0.2973032338908587, -0.0032350644379705897, 0, -0.000012947143763155565, -0.00006732564446879152, 0, -0.00033777251427091797, -0.013821083459497624, 0.0000000620253884096283, -0.00008689433851932771)
    5 = @(0.3192709837661525, 0.02669788538930942, 0.00025161284594377463, 0.0006999242398236953, 0.001984943605735288, -0.006396920237367126, -0.00005686705753037598, -0.0012302620062785738, 0, 0.000017433095657659603)
    6 = @(0.5526475353962192, 0.26199082542929136, 0, 0.00003698296957918854, 0.0002636603593251701, 0, -0.000014515526708548298, -0.02889974342231113, 0, -0.0000006581791093251965)
    7 = @(0.49122394864071234, 0, -0.0598994944903102, -0.0010324772290511621, -0.003621888338970474, 0.003005360656766159, 0, 0.0000648804178399764, 0, 0.00006003222821882037)
    8 = @(0.27679020651193353, -0.1795045635425241, 0, -0.00013294355343508696, -0.008966226743567987, 0, 0.0000219189880893702, -0.03321570046870835, 0, 0.007363773191367329)
    9 = @(0.11025296944415272, 0, -0.04043742101892383, -0.014145558952022252, -0.10626882866306611, -0.007531357804596345, -0.0013010031260203402, 0.002966066238598737, 0, 0.00018086625824931257)
    10 = @(0.19794785182838093, 0.11677208862318753, 0, -0.0004120552595670787, 0.004043398314355482, 0, -0.00014454481011724816, 0.03932842400353097, -0.08124058940123345, 0.00934816091407202)
    11 = @(0.3368479680992092, 0, 0.07938852346098085, 0.003095915891878116, 0.050313717646161535, -0.01578367133064049, -0.003851322103420169, 0.071226949947291, 0, -0.04548999724142255)
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($r, $i + 2).Value = $vals[$i]
    }
}

# --- Append the new row for the latest run (2025-08-30) ---------------------

$ws.Range("A12").NumberFormat = "@"
$ws.Range("A12").Value = "2025-08-30"
$ws.Range("A12").Style = "Normal"

$row12 = @(0.2151710127862062, -0.026704067178539403, 0, 0.0007020399580114204, 0.00001376601898416358, 0, 0.00003910068792868206, -0.05160562116775608, 0, -0.04412217363163179)

for ($i = 0; $i -lt $row12.Length; $i++) {
    $ws.Cells.Item(12, $i + 2).Value = $row12[$i]
}

# --- Column D needs to be a touch wider for the refreshed figures ----------
# (target stored width is 15.64453125; ColumnWidth is expressed in the
# "characters" unit Excel uses, which is the stored width minus the default
# 5/6 character margin)

$ws.Columns.Item(4).ColumnWidth = 14.811197916666666

